$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 19-32 (these rows no longer exist in the updated data)
$ws.Rows("19:32").Delete()

# Update the YEAR_WEEK (column B) values for rows 2-18 to reflect the
# corrected (non-COVID) week numbering for 2019
$newWeeks = @(201910,201911,201912,201913,201914,201915,201916,201917,201918,201919,201920,201921,201922,201923,201924,201925,201926)

for ($i = 0; $i -lt $newWeeks.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $newWeeks[$i]
}

# Rows 16-18 previously had YEAR = 2020; now they belong to 2019.
# Force text format so the value stays a string (matching column D's
# existing text-formatted entries) instead of being coerced to a number,
# then clear the format back off so no stray formatting lingers on the
# cell itself.
$yearRange = $ws.Range("D16:D18")
$yearRange.NumberFormat = "@"
$ws.Range("D16").Value = "2019"
$ws.Range("D17").Value = "2019"
$ws.Range("D18").Value = "2019"
$yearRange.ClearFormats()
